$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the two province rows (Huesca/Huelva) that got re-sorted, and
# carry their "Casos activos" (column C) values along with them.
$ws.Range("A53").Value = "Huelva"
$ws.Range("C53").Value = 72

$ws.Range("A54").Value = "Huesca"
$ws.Range("C54").Value = 0

# Update the "last updated" timestamp shown in A1.
$ws.Range("A1").Value = "Datos actualizados a 21 de Marzo de 2020 a las 23:46"
